# "updated chart of account output, placed a string newline replacer"
#
# The "Cost Center #" column (column B on the Data sheet) is removed
# entirely, so every column to its right (Journal No .. Date Deposited)
# shifts one slot to the left. The ChartofAccounts!I2:I3 list validation
# that used to live on column K follows the same shift and ends up on
# column J. The sheet's selection is also reset onto the new column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Re-home the "Journal Type" list validation (sourced from
# ChartofAccounts!$I$2:$I$3) onto a plain validation object covering its
# current range. Deleting column B below will shift this range left by
# one column automatically, landing it on J2:J1048576 -- matching where
# "Journal Type" (column K -> J) ends up.
$ws.Range("K2:K1048576").Validation.Add(3, 1, 1, "ChartofAccounts!`$I`$2:`$I`$3")

# Delete the "Cost Center #" column outright (column B), shifting every
# later column (Journal No, Account, Debit, Credit, Description, Name,
# Cheque No, Reference, Journal Type, Date Deposited) one place left.
$ws.Columns.Item(2).Delete()

# Match the saved selection state: active cell B1, whole column B
# selected, and no stale top-left scroll anchor.
$ws.Activate()
$ws.Range("B1:B1048576").Select()
